$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Convert B19 from a text/inline-string value to a true numeric value
$ws.Range("B19").Value = 23305567

# 2) Append new row 20 with the latest form submission
$ws.Range("A20").Value = "2025-10-26 16:03:14"

# B20 must stay a text value ("23201478"), not auto-convert to a number
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "23201478"
$ws.Range("B20").Style = "Normal"

$ws.Range("C20").Value = "sdsasda"
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = "CARLOS ALBERTO ORTEGA GALEANO"
$ws.Range("I20").Value = "CL 29 CR 50 D -18 (INTERIOR 1306 )"
$ws.Range("J20").Value = "Cumplido"
$ws.Range("K20").Value = "VENCIDO"
$ws.Range("L20").Value = 98535360
$ws.Range("M20").Value = "Formulario"
$ws.Range("N20").Value = 'https://drive.google.com/file/d/Error de carpeta: <HttpError 404 when requesting https://www.googleapis.com/drive/v3/files/1-Rg12PF0j59-sLYkjn3e_Hy9lcwxF5uz?fields=id%2C+name&alt=json returned "File not found: 1-Rg12PF0j59-sLYkjn3e_Hy9lcwxF5uz.". Details: "[{''message'': ''File not found: 1-Rg12PF0j59-sLYkjn3e_Hy9lcwxF5uz.'', ''domain'': ''global'', ''reason'': ''notFound'', ''location'': ''fileId'', ''locationType'': ''parameter''}]"</HttpError>/view, https://drive.google.com/file/d/Error de carpeta: <HttpError 404 when requesting https://www.googleapis.com/drive/v3/files/1-Rg12PF0j59-sLYkjn3e_Hy9lcwxF5uz?fields=id%2C+name&alt=json returned "File not found: 1-Rg12PF0j59-sLYkjn3e_Hy9lcwxF5uz.". Details: "[{''message'': ''File not found: 1-Rg12PF0j59-sLYkjn3e_Hy9lcwxF5uz.'', ''domain'': ''global'', ''reason'': ''notFound'', ''location'': ''fileId'', ''locationType'': ''parameter''}]"</HttpError>/view'
